# 1DES lima aula01 add
#
# - Moves the "Leiautes de tela" slide (with its screenshot picture) to the
#   end of the deck.
# - Inserts a brand new "Tags iniciais" slide (Title + Content layout) in
#   the slot the picture slide used to occupy, right before it, filled with
#   an explanation of the basic HTML tags.
# - The "Wireframe" slide that used to follow "Leiautes de tela" simply
#   slides up into position 8 unchanged.

$p = $ppt.ActivePresentation

# --- Step 1: move "Leiautes de tela" (+ screenshot) slide to the end ------
$leiautes = $p.Slides.Item(8)
$leiautes.MoveTo($p.Slides.Count)

# --- Step 2: insert a new Title+Content slide right before it ------------
$newIndex = $p.Slides.Count
$tagsSlide = $p.Slides.Add($newIndex, 2)

$title = $tagsSlide.Shapes.Item(1)
$title.Name = "Título 1"
$title.TextFrame.TextRange.Text = "Tags iniciais"

$body = $tagsSlide.Shapes.Item(2)
$body.Name = "Espaço Reservado para Conteúdo 2"

$lines = @(
    "<!DOCTYPE html> Define que o arquivo é uma página da  internet",
    "<HTML> Início da página – Tag principal",
    "<HEAD> Cabeçalho – Configurações da página",
    "<TITLE> Título que aparece na barra de títulos do navegador(Chrome)",
    "<META charset=`u{2018}utf-8`u{2019}> Configura os caracteres (Acentos em Português)",
    "<BODY> Corpo da página, onde colocamos tudo que será exibido ao usuário final (Front-End).",
    "<H1><H2>..<H5> Títulos e subtítulos da página",
    "<P> Paragrafo",
    "<center> alinha qualquer objeto ou texto",
    "<p align=`u{2018}center, rigth, left ou justfy`u{2019}> alinha o conteúdo do parágrafo.",
    "<font color size face> Altera a cor, o tamanho e o estilo da fonte",
    "<b><i><u> Negrito, Itálico e Sublinhado",
    "",
    "",
    ""
)
$body.TextFrame.TextRange.Text = [string]::Join([char]13, $lines)

# Ask the engine to shrink the text to fit the placeholder, same as
# PowerPoint's own "Shrink text on overflow" autofit would do once this much
# text is typed into the box (-> <a:normAutofit/> in the saved XML).
$body.TextFrame.AutoSize = 2
